# Weekly update: insert 3 new daily records (date 45244) for
# Arándano (blue) at Vega Central Mapocho de Santiago, ahead of the
# existing row that used to be row 304 (now shifted down to row 307).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 304; this shifts the
# old rows 304-367 down to 307-370 and already bumps the sheet
# dimension to A1:T370.
$ws.Rows("304:306").Insert()

# Columns that are constant across every record in this sheet.
$colA = 9
$colB = "Vega Central Mapocho de Santiago"
$colC = "Metropolitana"
$colE = 13
$colF = "Fruta"
$colG = 100101
$colH = "Berries"
$colI = 100101001
$colJ = "Arándano (blue)"
$colK = "Sin especificar"
$colQ = "`$/bandeja 2 kilos"
$colT = 2

function Set-Row($r, $d, $l, $m, $n, $o, $p, $rOrigin, $s) {
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $colF
    $ws.Cells.Item($r, 7).Value = $colG
    $ws.Cells.Item($r, 8).Value = $colH
    $ws.Cells.Item($r, 9).Value = $colI
    $ws.Cells.Item($r, 10).Value = $colJ
    $ws.Cells.Item($r, 11).Value = $colK
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $colQ
    $ws.Cells.Item($r, 18).Value = $rOrigin
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $colT
}

Set-Row 304 45244 "Especial" 210 7500 7500 7500 "Región del Maule" 3750
Set-Row 305 45244 "Primera"  280 6500 6500 6500 "Región del Maule" 3250
Set-Row 306 45244 "Segunda"  250 5000 5000 5000 "Región del Maule" 2500
